# New weekly price observations for "Piña" (Terminal Hortofrutícola Agro
# Chillán) are prepended to the data block: two new rows (Primera / Segunda
# quality, dated 2021-10-15) are inserted right before the current row 84,
# pushing every existing data row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 84; everything from the old row 84
# downward shifts to row 86 onward (formatting - e.g. the date style on
# column D - is carried over from the row being pushed down).
$ws.Rows("84:85").Insert()

# New row 84: Primera
$ws.Range("A84").Value = 7
$ws.Range("B84").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C84").Value = 'Ñuble'
$ws.Range("D84").Value = 44484
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = 'Fruta'
$ws.Range("G84").Value = 100108
$ws.Range("H84").Value = 'Tropicales y subtropicales'
$ws.Range("I84").Value = 100108005
$ws.Range("J84").Value = 'Piña'
$ws.Range("K84").Value = 'Caramelo'
$ws.Range("L84").Value = 'Primera'
$ws.Range("M84").Value = 60
$ws.Range("N84").Value = 18000
$ws.Range("O84").Value = 19000
$ws.Range("P84").Value = 18500
$ws.Range("Q84").Value = '$/caja 12 unidades'
$ws.Range("R84").Value = 'Ecuador'
$ws.Range("S84").Value = 1542
$ws.Range("T84").Value = 12

# New row 85: Segunda
$ws.Range("A85").Value = 7
$ws.Range("B85").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C85").Value = 'Ñuble'
$ws.Range("D85").Value = 44484
$ws.Range("E85").Value = 16
$ws.Range("F85").Value = 'Fruta'
$ws.Range("G85").Value = 100108
$ws.Range("H85").Value = 'Tropicales y subtropicales'
$ws.Range("I85").Value = 100108005
$ws.Range("J85").Value = 'Piña'
$ws.Range("K85").Value = 'Caramelo'
$ws.Range("L85").Value = 'Segunda'
$ws.Range("M85").Value = 60
$ws.Range("N85").Value = 18000
$ws.Range("O85").Value = 19000
$ws.Range("P85").Value = 18500
$ws.Range("Q85").Value = '$/caja 14 unidades'
$ws.Range("R85").Value = 'Ecuador'
$ws.Range("S85").Value = 1321
$ws.Range("T85").Value = 14
